$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2021-12-16T17:36:56+00:00"
$wsMeta.Range("B12").Value = "Status of the employee based on one or more code systems. Example codes include HIPAA (HipaaEmployeeStatusCodeSystem), Payer (WhPayerEmployeeStatusCodeSystem) or customer-specific codes."

# --- Sheet "Elements" updates ---
$wsEl = $wb.Worksheets.Item("Elements")

# Row 6 picks up the "closing" values that used to live on row 7's slice entry,
# and drops the slicing-specific values that no longer apply.
$wsEl.Range("AA6").Value = ""
$wsEl.Range("AB6").Value = ""
$wsEl.Range("AD6").Value = ""
$wsEl.Range("AI6").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`n"

# Row 7 (the valueCodeableConcept slice row) is removed entirely.
$wsEl.Rows.Item(7).Delete()

# Column width changes
$wsEl.Columns.Item(2).ColumnWidth = 12.65625
$wsEl.Columns.Item(25).ColumnWidth = 19.625
